$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2785.0334
$ws.Range("J17").Value = 2372.3215
$ws.Range("L17").Value = 7116.9645
$ws.Range("N17").Value = -7452.9645
$ws.Range("H53").Value = 2033.8
$ws.Range("I53").Value = 3350
$ws.Range("J53").Value = 59.5
$ws.Range("K53").Value = 3350
$ws.Range("L53").Value = 59.5
$ws.Range("M53").Value = -2713
$ws.Range("N53").Value = -1333.5
$ws.Range("H86").Value = 3160.6
$ws.Range("I86").Value = 3160.6
$ws.Range("K86").Value = 3160.6
$ws.Range("M86").Value = -2037.6
$ws.Range("H89").Value = 3160.6
$ws.Range("I89").Value = 3160.6
$ws.Range("K89").Value = 15803
$ws.Range("M89").Value = -10187
$ws.Range("H101").Value = 6899.3335
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H132").Value = 1371.1471
$ws.Range("I132").Value = 1269.375
$ws.Range("K132").Value = 3808.125
$ws.Range("M132").Value = -1278.125
$ws.Range("H138").Value = 4347.85
$ws.Range("I138").Value = 5124.3
$ws.Range("J138").Value = 3571.4
$ws.Range("K138").Value = 15372.9
$ws.Range("L138").Value = 10714.2
$ws.Range("M138").Value = -10232.9
$ws.Range("N138").Value = -20994.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1948.2858
$ws.Range("I2").Value = 1900
$ws.Range("J2").Value = 1967.6
$ws.Range("K2").Value = 1900
$ws.Range("L2").Value = 1967.6
$ws.Range("M2").Value = -1787
$ws.Range("N2").Value = -2193.6
$ws.Range("H32").Value = 4738.904
$ws.Range("I32").Value = 3753.9778
$ws.Range("K32").Value = 3753.9778
$ws.Range("M32").Value = -3466.9778
$ws.Range("H102").Value = 1693.25
$ws.Range("I102").Value = 1693.25
$ws.Range("K102").Value = 1693.25
$ws.Range("M102").Value = -71.25
$ws.Range("H116").Value = 1948.2858
$ws.Range("I116").Value = 1900
$ws.Range("J116").Value = 1967.6
$ws.Range("K116").Value = 1900
$ws.Range("L116").Value = 1967.6
$ws.Range("M116").Value = 394
$ws.Range("N116").Value = -6555.6
$ws.Range("H132").Value = 1380.0312
$ws.Range("I132").Value = 1116.5358
$ws.Range("K132").Value = 3349.6074
$ws.Range("M132").Value = -819.6074000000003

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1948.2858
$ws.Range("I3").Value = 1900
$ws.Range("J3").Value = 1967.6
$ws.Range("K3").Value = 1900
$ws.Range("L3").Value = 1967.6
$ws.Range("M3").Value = -1786
$ws.Range("N3").Value = -2195.6
$ws.Range("H134").Value = 6082.2085
$ws.Range("I134").Value = 6082.2085
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 18246.6255
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -15711.6255
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16250
$ws.Range("H58").Value = 3709.3157
$ws.Range("I58").Value = 2883.8572
$ws.Range("J58").Value = 4190.8335
$ws.Range("K58").Value = 2883.8572
$ws.Range("L58").Value = 4190.8335
$ws.Range("M58").Value = -2680.8572
$ws.Range("N58").Value = -4596.8335
$ws.Range("H62").Value = 2062.6667
$ws.Range("I62").Value = 2135
$ws.Range("K62").Value = 2135
$ws.Range("M62").Value = -1511
$ws.Range("H65").Value = 2062.6667
$ws.Range("I65").Value = 2135
$ws.Range("K65").Value = 10675
$ws.Range("M65").Value = -7555
$ws.Range("H132").Value = 2834.9
$ws.Range("I132").Value = 2169.8696
$ws.Range("K132").Value = 6509.6088
$ws.Range("M132").Value = -3979.6088
$ws.Range("H134").Value = 1908.88
$ws.Range("I134").Value = 1739.3334
$ws.Range("K134").Value = 5218.0002
$ws.Range("M134").Value = -2683.0002
$ws.Range("H136").Value = 3709.3157
$ws.Range("I136").Value = 2883.8572
$ws.Range("J136").Value = 4190.8335
$ws.Range("K136").Value = 8651.571599999999
$ws.Range("L136").Value = 12572.5005
$ws.Range("M136").Value = -6101.571599999999
$ws.Range("N136").Value = -17672.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 11745.777
$ws.Range("I34").Value = 20098
$ws.Range("J34").Value = 1305.5
$ws.Range("K34").Value = 60294
$ws.Range("L34").Value = 3916.5
$ws.Range("M34").Value = -60210
$ws.Range("N34").Value = -4084.5
$ws.Range("H121").Value = 921.4286
$ws.Range("I121").Value = 1030
$ws.Range("J121").Value = 903.3333
$ws.Range("K121").Value = 3090
$ws.Range("L121").Value = 2709.9999
$ws.Range("M121").Value = -1780
$ws.Range("N121").Value = -5329.9999
$ws.Range("H131").Value = 15296.2705
$ws.Range("J131").Value = 16621.842
$ws.Range("L131").Value = 49865.526
$ws.Range("N131").Value = -59945.526
$ws.Range("H140").Value = 2224.3704
$ws.Range("I140").Value = 1228.1578
$ws.Range("K140").Value = 3684.4734
$ws.Range("M140").Value = 1495.5266
$ws.Range("H141").Value = 2955.0908
$ws.Range("I141").Value = 3056.2222
$ws.Range("K141").Value = 9168.6666
$ws.Range("M141").Value = -3988.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 57011
$ws.Range("J27").Value = 57011
$ws.Range("L27").Value = 57011
$ws.Range("N27").Value = -57343
$ws.Range("H70").Value = 3499.8
$ws.Range("I70").Value = 2999.75
$ws.Range("K70").Value = 2999.75
$ws.Range("M70").Value = -2729.75
$ws.Range("H73").Value = 3499.8
$ws.Range("I73").Value = 2999.75
$ws.Range("K73").Value = 2999.75
$ws.Range("M73").Value = -2063.75
$ws.Range("H97").Value = 2106.9092
$ws.Range("J97").Value = 1974.25
$ws.Range("L97").Value = 1974.25
$ws.Range("N97").Value = -2966.25
$ws.Range("H102").Value = 3492.6667
$ws.Range("I102").Value = 3942.4736
$ws.Range("K102").Value = 3942.4736
$ws.Range("M102").Value = -2320.4736
$ws.Range("H132").Value = 2932.3157
$ws.Range("I132").Value = 2301.8
$ws.Range("J132").Value = 3632.889
$ws.Range("K132").Value = 6905.400000000001
$ws.Range("L132").Value = 10898.667
$ws.Range("M132").Value = -4375.400000000001
$ws.Range("N132").Value = -15958.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4848
$ws.Range("I7").Value = 3921
$ws.Range("J7").Value = 5775
$ws.Range("K7").Value = 3921
$ws.Range("L7").Value = 5775
$ws.Range("M7").Value = -3809
$ws.Range("N7").Value = -5999
$ws.Range("H55").Value = 362.12903
$ws.Range("I55").Value = 373.77777
$ws.Range("J55").Value = 346
$ws.Range("K55").Value = 373.77777
$ws.Range("L55").Value = 346
$ws.Range("M55").Value = -200.77777
$ws.Range("N55").Value = -692
$ws.Range("H105").Value = 15554.5
$ws.Range("J105").Value = 15554.5
$ws.Range("L105").Value = 15554.5
$ws.Range("N105").Value = -22542.5
$ws.Range("H126").Value = 4848
$ws.Range("I126").Value = 3921
$ws.Range("J126").Value = 5775
$ws.Range("K126").Value = 11763
$ws.Range("L126").Value = 17325
$ws.Range("M126").Value = -9293
$ws.Range("N126").Value = -22265
$ws.Range("H132").Value = 5463.125
$ws.Range("I132").Value = 5599.4
$ws.Range("K132").Value = 16798.2
$ws.Range("M132").Value = -14268.2
$ws.Range("H136").Value = 3712.0645
$ws.Range("I136").Value = 2927.0454
$ws.Range("J136").Value = 5631
$ws.Range("K136").Value = 8781.136200000001
$ws.Range("L136").Value = 16893
$ws.Range("M136").Value = -6231.136200000001
$ws.Range("N136").Value = -21993

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1244.7142
$ws.Range("I81").Value = 639.25
$ws.Range("J81").Value = 2052
$ws.Range("K81").Value = 1278.5
$ws.Range("L81").Value = 4104
$ws.Range("M81").Value = -217.5
$ws.Range("N81").Value = -6226
$ws.Range("H84").Value = 1244.7142
$ws.Range("I84").Value = 639.25
$ws.Range("J84").Value = 2052
$ws.Range("K84").Value = 6392.5
$ws.Range("L84").Value = 20520
$ws.Range("M84").Value = -1088.5
$ws.Range("N84").Value = -31128
$ws.Range("H100").Value = 609
$ws.Range("I100").Value = 510.5
$ws.Range("J100").Value = 1200
$ws.Range("K100").Value = 1021
$ws.Range("L100").Value = 2400
$ws.Range("M100").Value = -480
$ws.Range("N100").Value = -3482
$ws.Range("H105").Value = 50000
$ws.Range("J105").Value = 50000
$ws.Range("L105").Value = 50000
$ws.Range("N105").Value = -56988
$ws.Range("H107").Value = 1471.1428
$ws.Range("I107").Value = 697.5
$ws.Range("K107").Value = 2092.5
$ws.Range("M107").Value = -172.5
$ws.Range("H126").Value = 17335.625
$ws.Range("I126").Value = 22541
$ws.Range("K126").Value = 67623
$ws.Range("M126").Value = -65153
$ws.Range("H136").Value = 1368.2667
$ws.Range("I136").Value = 1339.7778
$ws.Range("J136").Value = 1411
$ws.Range("K136").Value = 4019.3334
$ws.Range("L136").Value = 4233
$ws.Range("M136").Value = -1469.3334
$ws.Range("N136").Value = -9333
